$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GanttChart")
$ws.Activate()

# --- Advance the "Display Week" scrollbar from week 4 to week 5 ---
# H4 is the cell linked to the Forms scrollbar control
# (xl/ctrlProps/ctrlProp1.xml, fmlaLink="$H$4"); every Week #/date header
# in the Gantt timeline recalculates off of it.
$ws.Range("H4").Value = 5

# Keep the scrollbar control itself in sync with its linked cell.
$scrollBar = $ws.Shapes.Item("Scroll Bar 46")
$scrollBar.ControlFormat.Value = 5

# --- Status update: mark the FAQ/Blog related tasks as 100% complete ---
$ws.Range("H25").Value = 1
$ws.Range("H26").Value = 1
$ws.Range("H27").Value = 1
$ws.Range("H28").Value = 1
$ws.Range("H29").Value = 1
$ws.Range("H30").Value = 1

# --- Update the saved view: scroll the frozen pane down and move the
#     active selection down to where the current work now is ---
$win = $excel.ActiveWindow
$win.ScrollRow = 22
[void]$ws.Range("R33").Select()
